$d = $word.ActiveDocument

function Find-ParagraphContaining($doc, $matchText) {
    foreach ($p in $doc.Paragraphs) {
        if ($p.Range.Text.Contains($matchText)) {
            return $p
        }
    }
    return $null
}

# ---------------------------------------------------------------------------
# 1) Insert a blank paragraph right after "All treatments thinning from below
#    to restore age structure with older trees" (before "Widespread reduction:
#    reduce everywhere to 35% of max biomass").
# ---------------------------------------------------------------------------
$pAllTreatments = Find-ParagraphContaining $d "All treatments thinning from below to restore age structure with older trees"
$pAllTreatments.Range.InsertParagraphAfter()

# ---------------------------------------------------------------------------
# 2) Insert a new paragraph with the "Patchy intensive treatment" text right
#    after the "Targeted reduction ... BAU otherwise" paragraph (before the
#    blank paragraph that precedes "Questions:").
# ---------------------------------------------------------------------------
$pBAU = Find-ParagraphContaining $d "no treatment in high carbon areas, BAU otherwise"
$pBAU.Range.InsertParagraphAfter()
$pPatchy = $pBAU.Next()
$pPatchy.Range.Text = "Patchy intensive treatment: many small treatments scattered across landscape, to slow down fire spread, vs "

# ---------------------------------------------------------------------------
# 3) Split the run "Biomass target not possible " into two runs:
#    "Residual b" + "iomass target not possible ", i.e. turn the sentence
#    into "Residual biomass target not possible ..." while leaving a run
#    boundary right after "Residual b".
# ---------------------------------------------------------------------------
$rngFind = $d.Content
$rngFind.Find.Execute("Biomass target not possible", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$leadingB = $d.Range($rngFind.Start, $rngFind.Start + 1)
$leadingB.Delete()
$insertionPoint = $d.Range($rngFind.Start, $rngFind.Start)
$insertionPoint.InsertBefore("Residual b")

# ---------------------------------------------------------------------------
# 4) The numbered (ListParagraph) bullet right after the "Biomass target not
#    possible" item is currently empty; give it the "Pace of treatments?"
#    text, then append two more numbered bullets after it.
# ---------------------------------------------------------------------------
$pEmptyBullet = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.ListFormat.ListType -ne 0 -and $p.Range.Text.Trim().Length -eq 0) {
        $pEmptyBullet = $p
        break
    }
}

$pEmptyBullet.Range.Text = "Pace of treatments?"

$pEmptyBullet.Range.InsertParagraphAfter()
$pOtherAreas = $pEmptyBullet.Next()
$pOtherAreas.Range.Text = "Other areas to avoid treating, aside from high carbon? Conservation areas or wilderness – let burn, follow other treatment plan, etc.?"

$pOtherAreas.Range.InsertParagraphAfter()
$pFireSuppression = $pOtherAreas.Next()
$pFireSuppression.Range.Text = "Fire suppression and Rx – high suppression at beginning, relaxing later once stands are restored?"

Write-Host "edit applied"
